# The workbook's cluster-outbreak table was replaced with an updated dataset
# (different cluster list, counts, and old/new flags). Since there is no row-for-row
# overlap between the old and new data, the simplest faithful reproduction is to
# clear the previous data rows and write the new dataset in full.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Old sheet had data through row 112; clear the whole former data range first so
# no stale rows linger past the new dataset (which ends at row 105).
$ws.Range("A2:C112").ClearContents()

# Cluster name / Active cases / Exist (old|new), one row per record.
$data = @(
    @('126 Racecourse Road Public Housing Tower Flemington', 5, 'old'),
    @('126 Racecourse Road Public Housing Tower Flemington', 7, 'new'),
    @('95 Napier Street Apartment Complex Fitzroy', 5, 'old'),
    @('Al Haj Halal Meats Glenroy', 11, 'new'),
    @('Al Haj Halal Meats Glenroy', 14, 'old'),
    @('Amiga Montessori Craigieburn', 7, 'old'),
    @('Budget Car and Truck Rentals Campbellfield', 5, 'old'),
    @('Campbellfield Ford Complex Vaccination Clinic Campbellfield', 8, 'new'),
    @('Campbellfield Ford Complex Vaccination Clinic Campbellfield', 9, 'old'),
    @('Cardinia Lakes Early Learning Centre Pakenham', 5, 'new'),
    @('Caroline Springs Police Station', 8, 'new'),
    @('Caroline Springs Police Station', 9, 'old'),
    @('Cedars Medical Clinic Coburg', 15, 'old'),
    @('Chemist Warehouse Campbellfield DC', 5, 'new'),
    @('Chemist Warehouse Campbellfield DC', 6, 'old'),
    @('Chemist Warehouse Fillo Drive Somerton', 11, 'old'),
    @('Chemist Warehouse Fillo Drive Somerton', 12, 'new'),
    @('City of Wyndham Community', 5, 'old'),
    @('City of Wyndham Community', 6, 'new'),
    @('Coles Campbellfield Plaza Campbellfield', 7, 'new'),
    @('Coles Campbellfield Plaza Campbellfield', 9, 'old'),
    @('Coles Coburg North Village', 17, 'new'),
    @('Coles Coburg North Village', 21, 'old'),
    @('Coles Pakenham Place Shopping Centre', 6, 'new'),
    @('Coles Pakenham Place Shopping Centre', 7, 'old'),
    @('Coles Roxburgh Village Roxburgh Park', 5, 'old'),
    @('Community Kids Bayswater Early Education Centre Bayswater North', 15, 'old'),
    @('Community Kids Bayswater Early Education Centre Bayswater North', 16, 'new'),
    @('Construction Site Olea Apartment Caulfield North', 14, 'new'),
    @('Construction Site Olea Apartment Caulfield North', 16, 'old'),
    @('Costco Wholesale Epping', 28, 'new'),
    @('Costco Wholesale Epping', 29, 'old'),
    @('Crusader Caravans Epping', 22, 'new'),
    @('Crusader Caravans Epping', 23, 'old'),
    @('Dandenong Police Station Dandenong', 5, 'old'),
    @('Dandenong Police Station Dandenong', 6, 'new'),
    @('Direct Freight Express Campbellfield', 7, 'new'),
    @('Direct Freight Express Campbellfield', 10, 'old'),
    @('Don Watson Coldstore Derrimut', 5, 'old'),
    @('Don Watson Coldstore Derrimutg', 5, 'new'),
    @('FedEx Station Melbourne Airport', 14, 'old'),
    @('Fine Food Holdings Pty Ltd Dandenong South', 9, 'old'),
    @('Fine Food Holdings Pty Ltd Dandenong South', 10, 'new'),
    @('Fitzroy Community School Fitzroy North', 5, 'old'),
    @('General Foods Campbellfield', 9, 'new'),
    @('Goodstart Early Learning Altona', 9, 'new'),
    @('Goodstart Early Learning Altona', 11, 'old'),
    @('Hamilton Marino 236 Jasper Road McKinnon', 9, 'new'),
    @('Hamilton Marino 236 Jasper Road McKinnon', 11, 'old'),
    @('Hello Fresh Warehouse Ravenhall', 5, 'old'),
    @('Hickory Construction Site Chadstone Car Park Malvern East', 5, 'old'),
    @('IGA Meadow Heights Shopping Centre Meadow Heights', 6, 'old'),
    @('ISS Factory Level 1 Terminal 2 Melbourne Airport Tullamarine', 10, 'new'),
    @('ISS Factory Level 1 Terminal 2 Melbourne Airport Tullamarine', 11, 'old'),
    @('Ibis Kingsgate Hotel Melbourne', 5, 'old'),
    @('Ibis Kingsgate Hotel Melbourne', 6, 'new'),
    @('Ilim College Kiewa Campus Boys Secondary Dallas', 5, 'new'),
    @('Ilim Learning Sanctuary Glenroy', 5, 'old'),
    @('Kasr Sweets Coolaroo', 5, 'old'),
    @('Kool Kidz Childcare Narre Warren', 12, 'old'),
    @('Kool Kidz Childcare Narre Warren', 13, 'new'),
    @('Lineage Logistics Laverton North', 6, 'new'),
    @('Lineage Logistics Laverton North', 7, 'old'),
    @('McDonald''s Craigieburn North', 5, 'old'),
    @('Mill Park Police Station Mill Park', 7, 'old'),
    @('Mill Park Police Station Mill Park', 8, 'new'),
    @('MyCentre Childcare Broadmeadows', 14, 'old'),
    @('MyCentre Childcare Broadmeadows', 21, 'new'),
    @('National Gallery of Victoria Melbourne', 8, 'new'),
    @('National Gallery of Victoria Melbourne', 9, 'old'),
    @('Nido Early School Ascot Vale', 11, 'old'),
    @('Nido Early School Ascot Vale', 12, 'new'),
    @('Nido Early School Glenroy', 19, 'old'),
    @('Nido Early School Glenroy', 22, 'new'),
    @('Northern Health Northern Hospital Epping Emergency Department Tier 1B', 55, 'new'),
    @('Northern Health Northern Hospital Epping Emergency Department Tier 1B', 66, 'old'),
    @('OnQ Plumbing and Excavations Craigieburn', 12, 'new'),
    @('OnQ Plumbing and Excavations Craigieburn', 13, 'old'),
    @('Oporto Coolaroo', 8, 'new'),
    @('Oporto Coolaroo', 9, 'old'),
    @('Our Lady Help of Christian''s Primary School Brunswick East', 9, 'new'),
    @('Our Lady Help of Christian''s Primary School Brunswick East', 10, 'old'),
    @('Panorama Construction Site Whitehorse Rd Box Hill', 5, 'old'),
    @('Private Residence Daycare Allumba Way Wollert', 8, 'new'),
    @('Ramsay Health Care Warrigal Private Hospital', 5, 'new'),
    @('Ramsay Health Care Warringal Private Hospital Heidelberg', 7, 'old'),
    @('St Vincents Hospital Emergency Department Melbourne', 8, 'old'),
    @('St Vincents Hospital Emergency Department Melbourne', 9, 'new'),
    @('Tek Foods Somerton', 10, 'new'),
    @('Tek Foods Somerton', 12, 'old'),
    @('The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1B', 18, 'new'),
    @('The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1B', 20, 'old'),
    @('The Royal Melbourne Hospital AMU Ward Parkville', 22, 'old'),
    @('The Royal Melbourne Hospital Ward 6SE Parkville', 21, 'new'),
    @('Truganina Early Learning Centre Truganina', 5, 'old'),
    @('Truganina Early Learning Centre Truganina', 6, 'new'),
    @('Unilodge College Square Student Accommodation 570 Lygon Street Carlton', 6, 'new'),
    @('Unilodge College Square Student Accommodation 570 Lygon Street Carlton', 7, 'old'),
    @('Wallaby Childcare Wollert', 17, 'old'),
    @('Wallaby Childcare Wollert', 18, 'new'),
    @('Western Health Footscray Hospital Emergency Department', 8, 'new'),
    @('Western Health Footscray Hospital Ward 3B Footscray', 10, 'old'),
    @('Yara Childcare Centre Truganina', 7, 'old'),
    @('Yarra Childcare Centre Truganina', 5, 'new')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
